$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "In the extraced folder ... dist folder, copy Cold Caller
#    Release.dmg file to Application folder" -- fix the "extraced" -> 
#    "extracted" typo (the run gets split the way Word's spell-check / 
#    AutoCorrect splits a corrected word, with the _GoBack bookmark ending
#    up between the "ct" and "ed folder" pieces) and mark "dist" and
#    "Release.dmg" as spell-check flagged words (w:proofErr).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(10)
$r1 = $p1.Range
Write-Output "p1 before: [$($r1.Text)]"

$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7E0AE5CE" w14:textId="42856A92" w:rsidR="009F6D52" w:rsidRDefault="00156B11" w:rsidP="00E55D2A"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>In the extra</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>ct</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>ed folder</w:t></w:r><w:r w:rsidR="009F6D52"><w:t xml:space="preserve"> &#8220;422P1-master&#8221;</w:t></w:r><w:r w:rsidR="009034FC"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00E55D2A"><w:t>go to</w:t></w:r><w:r w:rsidR="009034FC"><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> folder,  copy Cold Caller </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Release.dmg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> file to Application folder</w:t></w:r></w:p>
'@

$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) Paragraph "Double click the Cold Caller icon, run the application" --
#    drop the stray <w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr> that sits
#    inside the paragraph's <w:pPr> (paragraph mark run formatting).
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(12)
$r2 = $p2.Range
Write-Output "p2 before: [$($r2.Text)]"

$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="0AFE2159" w14:textId="20161CDC" w:rsidR="00E55D2A" w:rsidRDefault="00E55D2A" w:rsidP="00E55D2A"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>Double click the Cold Caller icon, run the application</w:t></w:r></w:p>
'@

$r2.InsertXML($xml2)

Write-Output "p1 after: [$($d.Paragraphs(10).Range.Text)]"
Write-Output "p2 after: [$($d.Paragraphs(12).Range.Text)]"
